$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (45406 -> 45436, one month later)
$ws.Range("A1").Value = 45436

# Update price list values in column D
$ws.Range("D28").Value = 230.1
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 336
$ws.Range("D31").Value = 422
